# Fixes issues #276 #262 #236 #193:
# Insert a new "B.1.17 - Dichiarazioni di uso standard beni culturali (BCS)"
# row right after the existing B.1.16 row (row 41), pushing the B.2 / C
# rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 42 (i.e. right after the current row 41,
# "B.1.16"); this shifts the old rows 42-51 down to 43-52 and copies the
# formatting of row 41 into the new row 42.
$ws.Rows.Item(42).Insert() | Out-Null

# Fill in the new row's values.
$ws.Range("A42").Value = "B"
$ws.Range("B42").Value = "Licenza Non Aperta"
$ws.Range("C42").Value = "B.1"
$ws.Range("D42").Value = "Solo uso non commerciale"
$ws.Range("E42").Value = "B.1.17"
$ws.Range("F42").Value = "Dichiarazioni di uso standard beni culturali (BCS)"

# Restore the view state (scroll position + active selection) to match the
# post-edit cursor position.
$ws.Range("F41").Select() | Out-Null
